$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new "Wins"/"Losses"/"Ties" columns, styled like the
# existing header cells (bold, centered, bordered) by copying AC1's format.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows (2-44): every player on the roster shares the team's overall
# 1999 record - 87 wins, 75 losses, 0 ties.
for ($r = 2; $r -le 44; $r++) {
    $ws.Cells.Item($r, 30).Value = 87   # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 75   # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF -> Ties
}
